$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.326.54"
$ws.Range("E2").Value = "  -0.90%  "
$ws.Range("D3").Value = "3.073.62"
$ws.Range("E3").Value = "  -1.56%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.64"
$ws.Range("E5").Value = "  -1.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.07"
$ws.Range("E6").Value = "  -1.42%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.071.18"
$ws.Range("E8").Value = "  -1.51%  "
$ws.Range("E9").Value = "  -2.18%  "
$ws.Range("E10").Value = "  -2.09%  "
$ws.Range("E11").Value = "  -2.53%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000239"
$ws.Range("E13").Value = "  -3.65%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.63"
$ws.Range("E14").Value = "  -4.37%  "
$ws.Range("E15").Value = "  -1.51%  "
$ws.Range("D16").Value = "3.585.33"
$ws.Range("E16").Value = "  -1.48%  "
$ws.Range("D17").Value = "66.236.02"
$ws.Range("E17").Value = "  -1.02%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.94"
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.59"
$ws.Range("E19").Value = "  +1.69%  "
$ws.Range("B20").Value = "WrappedEther"
$ws.Range("C20").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D20").Value = "3.070.00"
$ws.Range("E20").Value = "  -1.74%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "485.03"
$ws.Range("E21").Value = "  +1.98%  "
$ws.Range("B22").Value = "Polygon"
$ws.Range("C22").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.683"
$ws.Range("E22").Value = "  -3.67%  "
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.65"
$ws.Range("E23").Value = "  -2.72%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.24"
$ws.Range("E24").Value = "  -1.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.62"
$ws.Range("E25").Value = "  -4.60%  "
$ws.Range("E26").Value = "  -3.39%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.05"
$ws.Range("E27").Value = "  -2.75%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  -0.13%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.79"
$ws.Range("E29").Value = "  -1.75%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.24"
$ws.Range("E30").Value = "  -5.18%  "
$ws.Range("E31").Value = "  -3.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.61"
$ws.Range("E32").Value = "  -3.50%  "
$ws.Range("E33").Value = "  -3.19%  "
$ws.Range("E34").Value = "  -4.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.15%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "48.07"
$ws.Range("E36").Value = "  +2.48%  "
$ws.Range("E37").Value = "  -3.46%  "
$ws.Range("E38").Value = "  -5.04%  "
$ws.Range("E39").Value = "  -1.58%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.301"
$ws.Range("E40").Value = "  -3.82%  "
$ws.Range("E41").Value = "  -5.13%  "
$ws.Range("E42").Value = "  -4.55%  "
$ws.Range("D43").Value = "2.775.86"
$ws.Range("E43").Value = "  -1.71%  "
$ws.Range("E44").Value = "  -1.02%  "
$ws.Range("E45").Value = "  -2.76%  "
$ws.Range("B46").Value = "Bittensor"
$ws.Range("C46").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "365.27"
$ws.Range("E46").Value = "  -4.65%  "
$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "134.33"
$ws.Range("E47").Value = "  -1.26%  "
$ws.Range("E48").Value = "  +0.00%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "24.17"
$ws.Range("E49").Value = "  -3.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.15"
$ws.Range("E50").Value = "  -2.51%  "
$ws.Range("E51").Value = "  -2.25%  "
